$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original values for columns D, J, K, L, M, P (rows 2-40)
# before applying the weekly reshuffle of fruit/vegetable price records.
# Use .Value2 for reads (plain numeric/date-serial, avoids locale-formatted .Value).
$orig = @{}
for ($r = 2; $r -le 40; $r++) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

# Row permutation: new row -> source row (data pulled from source row's snapshot)
$map = @{
    2 = 11
    3 = 30
    4 = 16
    5 = 34
    6 = 35
    7 = 10
    8 = 31
    9 = 13
    10 = 40
    11 = 22
    12 = 8
    13 = 14
    14 = 29
    15 = 21
    16 = 12
    17 = 3
    18 = 6
    19 = 17
    20 = 32
    21 = 27
    22 = 33
    23 = 36
    24 = 4
    25 = 26
    26 = 39
    27 = 9
    28 = 19
    29 = 5
    30 = 20
    31 = 25
    32 = 23
    33 = 38
    34 = 15
    35 = 24
    36 = 28
    37 = 7
    38 = 2
    39 = 37
    40 = 18
}

foreach ($r in $map.Keys) {
    $src = $map[$r]
    $s = $orig[$src]
    $ws.Cells.Item($r, 4).Value = $s.D
    $ws.Cells.Item($r, 10).Value = $s.J
    $ws.Cells.Item($r, 11).Value = $s.K
    $ws.Cells.Item($r, 12).Value = $s.L
    $ws.Cells.Item($r, 13).Value = $s.M
    $ws.Cells.Item($r, 16).Value = $s.P
}

Write-Output "Reshuffled rows 2-40 (D,J,K,L,M,P) per weekly update."
